$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=id B=year C=week D=month E=day F=firm G=collectmode H=unsure I=n
#          J=c_poutou K=c_arthaud L=c_melenchon M=c_roussel N=c_montebourg
#          O=c_jadot P=c_hidalgo Q=c_macron R=c_pecresse S=c_barnier
#          T=c_bertrand U=c_lassalle V=c_daignant W=c_lepen X=c_zemmour

# New row 53 - ifop poll (10/11)
$ws.Cells.Item(53, 1).Value = 15
$ws.Cells.Item(53, 2).Value = 2021
$ws.Cells.Item(53, 3).Value = 6
$ws.Cells.Item(53, 4).Value = 10
$ws.Cells.Item(53, 5).Value = 6
$ws.Cells.Item(53, 6).Value = "ifop"
$ws.Cells.Item(53, 7).Value = "online"
$ws.Cells.Item(53, 8).Value = "included"
$ws.Cells.Item(53, 9).Value = 921
$ws.Cells.Item(53, 10).Value = 0.5
$ws.Cells.Item(53, 11).Value = 0.5
$ws.Cells.Item(53, 12).Value = 7
$ws.Cells.Item(53, 13).Value = 2
$ws.Cells.Item(53, 14).Value = 2
$ws.Cells.Item(53, 15).Value = 8
$ws.Cells.Item(53, 16).Value = 6
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 20).Value = 16
$ws.Cells.Item(53, 21).Value = 1
$ws.Cells.Item(53, 22).Value = 2
$ws.Cells.Item(53, 23).Value = 16
$ws.Cells.Item(53, 24).Value = 14

# New row 54 - ifop poll (10/11)
$ws.Cells.Item(54, 1).Value = 15
$ws.Cells.Item(54, 2).Value = 2021
$ws.Cells.Item(54, 3).Value = 6
$ws.Cells.Item(54, 4).Value = 10
$ws.Cells.Item(54, 5).Value = 6
$ws.Cells.Item(54, 6).Value = "ifop"
$ws.Cells.Item(54, 7).Value = "online"
$ws.Cells.Item(54, 8).Value = "included"
$ws.Cells.Item(54, 9).Value = 921
$ws.Cells.Item(54, 10).Value = 0.5
$ws.Cells.Item(54, 11).Value = 0.5
$ws.Cells.Item(54, 12).Value = 8
$ws.Cells.Item(54, 13).Value = 1.5
$ws.Cells.Item(54, 14).Value = 2.5
$ws.Cells.Item(54, 15).Value = 8
$ws.Cells.Item(54, 16).Value = 6
$ws.Cells.Item(54, 17).Value = 26
$ws.Cells.Item(54, 18).Value = 11
$ws.Cells.Item(54, 21).Value = 1
$ws.Cells.Item(54, 22).Value = 3
$ws.Cells.Item(54, 23).Value = 18
$ws.Cells.Item(54, 24).Value = 14

# Match the author's final selection / scroll position
$ws.Range("V54").Select()
